# The deck's single design ("Integral" / "Red Violet" colour scheme, stored
# in ppt/theme/theme1.xml) is switched over to the stock "Office" colour
# palette. Every slide shares the one Slide Master/theme in this file, so
# grabbing the ThemeColorScheme from any slide reaches the same underlying
# <a:clrScheme> - walk all twelve theme colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) and set each to the matching "Office" RGB
# value via the real PowerPoint automation surface
# (ThemeColorScheme.Colors(i).RGB), exactly as a user would if they picked
# the built-in "Office" colour variant for the current design.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (slot name, new "Office" RGB as a COM long 0xBBGGRR)
$officeColors = @(
    0,        # 1  dk1      000000
    16777215, # 2  lt1      FFFFFF
    6968388,  # 3  dk2      44546A
    15132391, # 4  lt2      E7E6E6
    13998939, # 5  accent1  5B9BD5
    3243501,  # 6  accent2  ED7D31
    10855845, # 7  accent3  A5A5A5
    49407,    # 8  accent4  FFC000
    12874308, # 9  accent5  4472C4
    4697456,  # 10 accent6  70AD47
    12673797, # 11 hlink    0563C1
    7491477   # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
